# Updates the crypto price/volume table (and the OKB/dogwifhat row order)
# to match the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophes force values that look numeric (e.g. "1.00") to stay
# stored as text, matching the original inline-string cell contents.
$ws.Range("D2").Value = "74.942.49"
$ws.Range("E2").Value = "  +7.14%  "
$ws.Range("D3").Value = "2.663.96"
$ws.Range("E3").Value = "  +8.63%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'186.36"
$ws.Range("E5").Value = "  +11.82%  "
$ws.Range("D6").Value = "'586.57"
$ws.Range("E6").Value = "  +2.89%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.534"
$ws.Range("E8").Value = "  +3.90%  "
$ws.Range("D9").Value = "'0.194"
$ws.Range("E9").Value = "  +10.66%  "
$ws.Range("D10").Value = "2.663.00"
$ws.Range("E10").Value = "  +8.66%  "
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("E12").Value = "  +5.66%  "
$ws.Range("D13").Value = "'4.73"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "3.153.76"
$ws.Range("E14").Value = "  +8.78%  "
$ws.Range("D15").Value = "74.534.86"
$ws.Range("E15").Value = "  +6.71%  "
$ws.Range("D16").Value = "'0.0000186"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "'26.49"
$ws.Range("E17").Value = "  +9.32%  "
$ws.Range("D18").Value = "2.661.66"
$ws.Range("E18").Value = "  +8.56%  "
$ws.Range("D19").Value = "'9.19"
$ws.Range("E19").Value = "  +28.15%  "
$ws.Range("E20").Value = "  +8.80%  "
$ws.Range("D21").Value = "'370.72"
$ws.Range("E21").Value = "  +8.50%  "
$ws.Range("D22").Value = "'2.26"
$ws.Range("E22").Value = "  +11.62%  "
$ws.Range("D23").Value = "'4.07"
$ws.Range("E23").Value = "  +4.35%  "
$ws.Range("D24").Value = "'6.25"
$ws.Range("E24").Value = "  +3.66%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'69.47"
$ws.Range("E26").Value = "  +4.59%  "
$ws.Range("D27").Value = "'4.14"
$ws.Range("E27").Value = "  +7.60%  "
$ws.Range("D28").Value = "'9.33"
$ws.Range("E28").Value = "  +9.45%  "
$ws.Range("D29").Value = "2.798.95"
$ws.Range("E29").Value = "  +8.53%  "
$ws.Range("D30").Value = "'1.01"
$ws.Range("E30").Value = "  +11.07%  "
$ws.Range("D31").Value = "0.0₃0943"
$ws.Range("E31").Value = "  +9.85%  "
$ws.Range("E32").Value = "  +13.22%  "
$ws.Range("D33").Value = "'522.14"
$ws.Range("E33").Value = "  +13.31%  "
$ws.Range("D34").Value = "'7.60"
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("D35").Value = "'1.75"
$ws.Range("E35").Value = "  +7.15%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'163.43"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("E38").Value = "  +5.28%  "
$ws.Range("D39").Value = "'19.19"
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'170.06"
$ws.Range("E42").Value = "  +26.61%  "
$ws.Range("D43").Value = "'4.99"
$ws.Range("E43").Value = "  +12.55%  "
$ws.Range("D44").Value = "'0.328"
$ws.Range("E44").Value = "  +7.86%  "
$ws.Range("D45").Value = "'1.67"
$ws.Range("E45").Value = "  +8.82%  "
$ws.Range("E46").Value = "  +8.22%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.35"
$ws.Range("E47").Value = "  +10.12%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'38.99"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("E49").Value = "  +15.62%  "
$ws.Range("E50").Value = "  +6.78%  "
$ws.Range("D51").Value = "'21.26"
$ws.Range("E51").Value = "  +21.47%  "
